$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 22
$ws.Range("B2").Value = 232

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 175

$ws.Range("B4").Value = 165

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 58
